$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Some Price values look numeric (e.g. "226.16"); force text format so
# Excel keeps them as literal text instead of auto-converting to a number,
# matching the original text-cell representation.

$ws.Range("D2").Value = "34.148.43"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.16"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.65"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0686"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "1.795.65"
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.98"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").Value = "34.054.33"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.20"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.02"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +4.46%  "
$ws.Range("D20").Value = "0.0₃0776"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.88"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +4.08%  "
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.37"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  +3.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.29"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0518"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.67"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.62"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +4.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.79"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("D35").Value = "1.444.44"
$ws.Range("E35").Value = "  +5.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.653"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("E37").Value = "  +10.38%  "
$ws.Range("E38").Value = "  +4.64%  "
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.24"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +4.38%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("E42").Value = "  +3.34%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.47"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("E45").Value = "  +5.04%  "
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "1.944.40"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.83"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("E51").Value = "  -0.03%  "
